$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last two data rows (old rows 4 and 5); row 3 takes over with
# updated values, and row 2 is updated in place.
$ws.Rows("4:5").Delete()

# Row 2: FAPs | Ccl21b | Cxcr3 | Resolving-Mac (D2 "ECs" -> "Resolving-Mac")
$ws.Cells.Item(2, 4).Value = "Resolving-Mac"
$ws.Cells.Item(2, 8).Value = 0.6219589999999999
$ws.Cells.Item(2, 9).Value = 0.6398583988494134
$ws.Cells.Item(2, 10).Value = 0.6398583988494134
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.888791333333333
$ws.Cells.Item(2, 14).Value = 5.666374
$ws.Cells.Item(2, 15).Value = 1
$ws.Cells.Item(2, 16).Value = 1
$ws.Cells.Item(2, 17).Value = 0.3915835896295556
$ws.Cells.Item(2, 18).Value = 3.524252306666
$ws.Cells.Item(2, 19).Value = 0.6398583988494134
$ws.Cells.Item(2, 20).Value = 0.6398583988494134

# Row 3: MuSCs | Ccl21b | Cxcr3 | Resolving-Mac
$ws.Cells.Item(3, 1).Value = "MuSCs"
$ws.Cells.Item(3, 4).Value = "Resolving-Mac"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.116689
$ws.Cells.Item(3, 8).Value = 0.350067
$ws.Cells.Item(3, 9).Value = 0.3601416011505865
$ws.Cells.Item(3, 10).Value = 0.3601416011505865
$ws.Cells.Item(3, 13).Value = 1.888791333333333
$ws.Cells.Item(3, 14).Value = 5.666374
$ws.Cells.Item(3, 15).Value = 1
$ws.Cells.Item(3, 16).Value = 1
$ws.Cells.Item(3, 17).Value = 0.2204011718953333
$ws.Cells.Item(3, 18).Value = 1.983610547058
$ws.Cells.Item(3, 19).Value = 0.3601416011505865
$ws.Cells.Item(3, 20).Value = 0.3601416011505865
